$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B86 was stored as text "3"; convert it to a real number 3
$ws.Cells.Item(86, 2).Value = 3

# Add new row 87 with the additional annotation data
$ws.Cells.Item(87, 1).Value = "Ying Tang"

# B87 must remain a text value "4" (not a number), so force text with a
# leading apostrophe and then reset the style so no extra formatting sticks
$ws.Cells.Item(87, 2).Value = "'4"
$ws.Cells.Item(87, 2).Style = "Normal"

$ws.Cells.Item(87, 3).Value = " I don't think ,I'm concerned that"
$ws.Cells.Item(87, 4).Value = "DFT"
$ws.Cells.Item(87, 5).Value = "WRI"
$ws.Cells.Item(87, 6).Value = "4d925051-98cd-4c47-ad57-530c84f23548"
$ws.Cells.Item(87, 7).Value = "B1QRgziT-_annotated.xlsx"
$ws.Cells.Item(87, 8).Value = "I don't think this paper explains the importance of its results nearly enough and I'm concerned that it may not be obvious what a breakthrough it is just from skimming the abstract."
